# Boolean Peaking and Flexibility Flags.xlsx -- update to add new fuel-type
# rows (crude oil, heavy or residual fuel oil, municipal solid waste) to both
# flag sheets, flip the "hard coal" flag off on BPaFF-BITPTaP, right-align
# the "Boolean" header cell on both flag sheets, and make "About" the active
# tab again (instead of BPaFF-BITPTaP).

$wb = $excel.ActiveWorkbook

$wsAbout   = $wb.Worksheets.Item("About")
$wsPeaker  = $wb.Worksheets.Item("BPaFF-BITPTaP")
$wsFlex    = $wb.Worksheets.Item("BPaFF-BDTPTPF")

# --- BPaFF-BITPTaP ("Is This Plant Type a Peaker") ---------------------
# hard coal flag flips from 1 to 0
$wsPeaker.Range("B2").Value = 0

# header cell (B1, "Boolean") becomes right aligned
$wsPeaker.Range("B1").HorizontalAlignment = -4152   # xlRight

# new fuel-type rows
$wsPeaker.Range("A15").Value = "crude oil"
$wsPeaker.Range("B15").Formula = "=B11"

$wsPeaker.Range("A16").Value = "heavy or residual fuel oil"
$wsPeaker.Range("B16").Formula = "=B11"

$wsPeaker.Range("A17").Value = "municipal solid waste"
$wsPeaker.Range("B17").Formula = "=B9"

# --- BPaFF-BDTPTPF ("Does This Plant Type Provide Flexibility") --------
# header cell (B1, "Boolean") becomes right aligned
$wsFlex.Range("B1").HorizontalAlignment = -4152   # xlRight

# new fuel-type rows
$wsFlex.Range("A15").Value = "crude oil"
$wsFlex.Range("B15").Formula = "=B11"

$wsFlex.Range("A16").Value = "heavy or residual fuel oil"
$wsFlex.Range("B16").Formula = "=B11"

$wsFlex.Range("A17").Value = "municipal solid waste"
$wsFlex.Range("B17").Formula = "=B9"

# --- Active tab goes back to "About" ------------------------------------
$wsAbout.Activate()
